$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 242, shifting existing rows 242-318 down to 243-319.
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new data record
# (same fixed market/category fields as surrounding rows; new date/price figures).
$ws.Cells.Item(242, 1).Value = 4
$ws.Cells.Item(242, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(242, 3).Value = "Los Lagos"
$ws.Cells.Item(242, 4).Value = 44809
$ws.Cells.Item(242, 5).Value = 10
$ws.Cells.Item(242, 6).Value = 100112043
$ws.Cells.Item(242, 7).Value = "Pepino ensalada"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 70
$ws.Cells.Item(242, 11).Value = 30000
$ws.Cells.Item(242, 12).Value = 30000
$ws.Cells.Item(242, 13).Value = 30000
$ws.Cells.Item(242, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(242, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(242, 16).Value = 500
$ws.Cells.Item(242, 17).Value = 60
$ws.Cells.Item(242, 18).Value = "Hortaliza"
